# "add a rough richfield" - fill in rough/preliminary habitat assessment
# values (hab_value, and for the first Richfield Creek crossing also the
# upstream habitat length + species codes) that were left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - reference_number 12, Richfield Creek (197663_us)
$ws.Range("H13").Value = "high"
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = "CH,CO,LNC,LSU,RB,ST"

# Row 14 - reference_number 13, Richfield Creek (197663_ds)
$ws.Range("H14").Value = "high"

# Row 35
$ws.Range("H35").Value = "moderate"

# Row 37
$ws.Range("H37").Value = "moderate"

# Row 38
$ws.Range("H38").Value = "moderate"

# Row 39
$ws.Range("H39").Value = "moderate"

# Leave the selection where the editor ended up after the last entry.
$ws.Range("M17").Select() | Out-Null
